$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values in rows 2-5 (rows shifted earlier than row 6 which is removed) ---
$ws.Cells.Item(2, 1).Value = 45057.50694444445
$ws.Cells.Item(2, 2).Value = 17.295
$ws.Cells.Item(2, 3).Value = 11.407
$ws.Cells.Item(2, 4).Value = 3.988
$ws.Cells.Item(2, 5).Value = 36.822
$ws.Cells.Item(2, 6).Value = 28.934
$ws.Cells.Item(2, 7).Value = 13.611
$ws.Cells.Item(2, 8).Value = 42.104
$ws.Cells.Item(2, 9).Value = 20.942
$ws.Cells.Item(2, 10).Value = 8.632
$ws.Cells.Item(2, 11).Value = 12.848
$ws.Cells.Item(2, 12).Value = 14.482
$ws.Cells.Item(2, 13).Value = 15.027
$ws.Cells.Item(2, 14).Value = 4.344
$ws.Cells.Item(2, 15).Value = 13.535
$ws.Cells.Item(2, 16).Value = 18.799
$ws.Cells.Item(2, 17).Value = 11.717
$ws.Cells.Item(2, 18).Value = 3.386
$ws.Cells.Item(2, 19).Value = 2.226
$ws.Cells.Item(2, 20).Value = 198.649
$ws.Cells.Item(2, 21).Value = 37.584
$ws.Cells.Item(2, 22).Value = 12.493
$ws.Cells.Item(2, 23).Value = 24.537
$ws.Cells.Item(2, 24).Value = 12.435
$ws.Cells.Item(2, 25).Value = 3.148
$ws.Cells.Item(2, 26).Value = 21.511
$ws.Cells.Item(2, 27).Value = 11.035
$ws.Cells.Item(2, 28).Value = 10.064
$ws.Cells.Item(2, 29).Value = 11.833
$ws.Cells.Item(2, 30).Value = 15.025
$ws.Cells.Item(2, 31).Value = 3.317
$ws.Cells.Item(2, 32).Value = 37.614
$ws.Cells.Item(2, 33).Value = 6.758
$ws.Cells.Item(2, 34).Value = 15.619

$ws.Cells.Item(3, 1).Value = 45057.51388888889
$ws.Cells.Item(3, 2).Value = 14.893
$ws.Cells.Item(3, 3).Value = 10.464
$ws.Cells.Item(3, 4).Value = 1.76
$ws.Cells.Item(3, 5).Value = 32.234
$ws.Cells.Item(3, 6).Value = 25.733
$ws.Cells.Item(3, 7).Value = 11.721
$ws.Cells.Item(3, 8).Value = 45.119
$ws.Cells.Item(3, 9).Value = 18.033
$ws.Cells.Item(3, 10).Value = 7.76
$ws.Cells.Item(3, 11).Value = 11.353
$ws.Cells.Item(3, 12).Value = 12.889
$ws.Cells.Item(3, 13).Value = 13.464
$ws.Cells.Item(3, 14).Value = 3.744
$ws.Cells.Item(3, 15).Value = 11.655
$ws.Cells.Item(3, 16).Value = 16.395
$ws.Cells.Item(3, 17).Value = 10.104
$ws.Cells.Item(3, 18).Value = 1.517
$ws.Cells.Item(3, 19).Value = 1.04
$ws.Cells.Item(3, 20).Value = 170.081
$ws.Cells.Item(3, 21).Value = 32.624
$ws.Cells.Item(3, 22).Value = 10.758
$ws.Cells.Item(3, 23).Value = 21.557
$ws.Cells.Item(3, 24).Value = 11.176
$ws.Cells.Item(3, 25).Value = 2.195
$ws.Cells.Item(3, 26).Value = 21.999
$ws.Cells.Item(3, 27).Value = 9.502000000000001
$ws.Cells.Item(3, 28).Value = 8.617000000000001
$ws.Cells.Item(3, 29).Value = 10.108
$ws.Cells.Item(3, 30).Value = 13.426
$ws.Cells.Item(3, 31).Value = 1.247
$ws.Cells.Item(3, 32).Value = 41.259
$ws.Cells.Item(3, 33).Value = 5.899
$ws.Cells.Item(3, 34).Value = 13.45

$ws.Cells.Item(4, 1).Value = 45057.52083333334
$ws.Cells.Item(4, 2).Value = 13.452
$ws.Cells.Item(4, 3).Value = 9.638
$ws.Cells.Item(4, 4).Value = 1.223
$ws.Cells.Item(4, 5).Value = 29.193
$ws.Cells.Item(4, 6).Value = 23.465
$ws.Cells.Item(4, 7).Value = 10.587
$ws.Cells.Item(4, 8).Value = 41.745
$ws.Cells.Item(4, 9).Value = 16.288
$ws.Cells.Item(4, 10).Value = 7.097
$ws.Cells.Item(4, 11).Value = 10.387
$ws.Cells.Item(4, 12).Value = 11.705
$ws.Cells.Item(4, 13).Value = 12.273
$ws.Cells.Item(4, 14).Value = 3.382
$ws.Cells.Item(4, 15).Value = 10.527
$ws.Cells.Item(4, 16).Value = 14.868
$ws.Cells.Item(4, 17).Value = 9.064
$ws.Cells.Item(4, 18).Value = 1.012
$ws.Cells.Item(4, 19).Value = 0.743
$ws.Cells.Item(4, 20).Value = 152.91
$ws.Cells.Item(4, 21).Value = 29.466
$ws.Cells.Item(4, 22).Value = 9.717000000000001
$ws.Cells.Item(4, 23).Value = 19.578
$ws.Cells.Item(4, 24).Value = 10.211
$ws.Cells.Item(4, 25).Value = 1.822
$ws.Cells.Item(4, 26).Value = 20.106
$ws.Cells.Item(4, 27).Value = 8.583
$ws.Cells.Item(4, 28).Value = 7.729
$ws.Cells.Item(4, 29).Value = 9.069000000000001
$ws.Cells.Item(4, 30).Value = 12.231
$ws.Cells.Item(4, 31).Value = 0.766
$ws.Cells.Item(4, 32).Value = 38.034
$ws.Cells.Item(4, 33).Value = 5.368
$ws.Cells.Item(4, 34).Value = 12.148

$ws.Cells.Item(5, 1).Value = 45057.52777777778
$ws.Cells.Item(5, 2).Value = 2.4
$ws.Cells.Item(5, 3).Value = 1.46
$ws.Cells.Item(5, 4).Value = 0.64
$ws.Cells.Item(5, 5).Value = 5.19
$ws.Cells.Item(5, 6).Value = 3.76
$ws.Cells.Item(5, 7).Value = 1.89
$ws.Cells.Item(5, 8).Value = 12.94
$ws.Cells.Item(5, 9).Value = 2.91
$ws.Cells.Item(5, 10).Value = 1.24
$ws.Cells.Item(5, 11).Value = 1.53
$ws.Cells.Item(5, 12).Value = 2.08
$ws.Cells.Item(5, 13).Value = 2.12
$ws.Cells.Item(5, 14).Value = 0.62
$ws.Cells.Item(5, 15).Value = 1.88
$ws.Cells.Item(5, 16).Value = 2.68
$ws.Cells.Item(5, 17).Value = 1.83
$ws.Cells.Item(5, 18).Value = 0.6899999999999999
$ws.Cells.Item(5, 19).Value = 0.31
$ws.Cells.Item(5, 20).Value = 21.38
$ws.Cells.Item(5, 21).Value = 5.64
$ws.Cells.Item(5, 22).Value = 1.74
$ws.Cells.Item(5, 23).Value = 3.64
$ws.Cells.Item(5, 24).Value = 1.77
$ws.Cells.Item(5, 25).Value = 0.58
$ws.Cells.Item(5, 26).Value = 5.78
$ws.Cells.Item(5, 27).Value = 1.53
$ws.Cells.Item(5, 28).Value = 1.51
$ws.Cells.Item(5, 29).Value = 1.74
$ws.Cells.Item(5, 30).Value = 2.14
$ws.Cells.Item(5, 31).Value = 0.55
$ws.Cells.Item(5, 32).Value = 12.3
$ws.Cells.Item(5, 33).Value = 0.85
$ws.Cells.Item(5, 34).Value = 2.18

# --- Remove row 6 (dataset shrank from 5 data rows to 4) ---
$ws.Rows.Item(6).Delete()

# --- Widen specific columns from 7 to 8 characters (stored OOXML width units) ---
$wideCols = 2,3,7,11,12,13,15,17,22,24,27,28,29,30,34
foreach ($c in $wideCols) {
    $ws.Columns.Item($c).ColumnWidth = 7.166666666666667
}
